$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (ownTeam, oppTeam) after column C, shifting old D:I to F:K
$ws.Columns("D:E").Insert()

# Mark columns A:K as text-formatted so numeric-looking strings are stored as text (matches source data)
$ws.Range("A1:K12").NumberFormat = "@"

# Header row
$ws.Cells.Item(1, 1).Value = "venue"
$ws.Cells.Item(1, 2).Value = "date"
$ws.Cells.Item(1, 3).Value = "result"
$ws.Cells.Item(1, 4).Value = "ownTeam"
$ws.Cells.Item(1, 5).Value = "oppTeam"
$ws.Cells.Item(1, 6).Value = "batsman"
$ws.Cells.Item(1, 7).Value = "totalRuns"
$ws.Cells.Item(1, 8).Value = "totalBalls"
$ws.Cells.Item(1, 9).Value = "total4s"
$ws.Cells.Item(1, 10).Value = "total6s"
$ws.Cells.Item(1, 11).Value = "sr"

# Data rows
# Row 2
$ws.Cells.Item(2, 1).Value = " Abu Dhabi"
$ws.Cells.Item(2, 2).Value = " October 07 2020"
$ws.Cells.Item(2, 3).Value = "KKR won by 10 runs"
$ws.Cells.Item(2, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(2, 5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(2, 6).Value = "Ravindra Jadeja "
$ws.Cells.Item(2, 7).Value = "21"
$ws.Cells.Item(2, 8).Value = "8"
$ws.Cells.Item(2, 9).Value = "3"
$ws.Cells.Item(2, 10).Value = "1"
$ws.Cells.Item(2, 11).Value = "262.50"

# Row 3
$ws.Cells.Item(3, 1).Value = " Sharjah"
$ws.Cells.Item(3, 2).Value = " October 23 2020"
$ws.Cells.Item(3, 3).Value = "Mumbai won by 10 wickets (with 46 balls remaining)"
$ws.Cells.Item(3, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(3, 5).Value = "Mumbai Indians"
$ws.Cells.Item(3, 6).Value = "Ravindra Jadeja "
$ws.Cells.Item(3, 7).Value = "7"
$ws.Cells.Item(3, 8).Value = "6"
$ws.Cells.Item(3, 9).Value = "1"
$ws.Cells.Item(3, 10).Value = "0"
$ws.Cells.Item(3, 11).Value = "116.66"

# Row 4
$ws.Cells.Item(4, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(4, 2).Value = " October 13 2020"
$ws.Cells.Item(4, 3).Value = "Super Kings won by 20 runs"
$ws.Cells.Item(4, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(4, 5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(4, 6).Value = "Ravindra Jadeja "
$ws.Cells.Item(4, 7).Value = "25"
$ws.Cells.Item(4, 8).Value = "10"
$ws.Cells.Item(4, 9).Value = "3"
$ws.Cells.Item(4, 10).Value = "1"
$ws.Cells.Item(4, 11).Value = "250.00"

# Row 5
$ws.Cells.Item(5, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(5, 2).Value = " September 25 2020"
$ws.Cells.Item(5, 3).Value = "Capitals won by 44 runs"
$ws.Cells.Item(5, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(5, 5).Value = "Delhi Capitals"
$ws.Cells.Item(5, 6).Value = "Ravindra Jadeja "
$ws.Cells.Item(5, 7).Value = "12"
$ws.Cells.Item(5, 8).Value = "9"
$ws.Cells.Item(5, 9).Value = "1"
$ws.Cells.Item(5, 10).Value = "0"
$ws.Cells.Item(5, 11).Value = "133.33"

# Row 6
$ws.Cells.Item(6, 1).Value = " Sharjah"
$ws.Cells.Item(6, 2).Value = " October 17 2020"
$ws.Cells.Item(6, 3).Value = "Capitals won by 5 wickets (with 1 ball remaining)"
$ws.Cells.Item(6, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(6, 5).Value = "Delhi Capitals"
$ws.Cells.Item(6, 6).Value = "Ravindra Jadeja "
$ws.Cells.Item(6, 7).Value = "33"
$ws.Cells.Item(6, 8).Value = "13"
$ws.Cells.Item(6, 9).Value = "0"
$ws.Cells.Item(6, 10).Value = "4"
$ws.Cells.Item(6, 11).Value = "253.84"

# Row 7
$ws.Cells.Item(7, 1).Value = " Abu Dhabi"
$ws.Cells.Item(7, 2).Value = " September 19 2020"
$ws.Cells.Item(7, 3).Value = "Super Kings won by 5 wickets (with 4 balls remaining)"
$ws.Cells.Item(7, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(7, 5).Value = "Mumbai Indians"
$ws.Cells.Item(7, 6).Value = "Ravindra Jadeja "
$ws.Cells.Item(7, 7).Value = "10"
$ws.Cells.Item(7, 8).Value = "5"
$ws.Cells.Item(7, 9).Value = "2"
$ws.Cells.Item(7, 10).Value = "0"
$ws.Cells.Item(7, 11).Value = "200.00"

# Row 8
$ws.Cells.Item(8, 1).Value = " Sharjah"
$ws.Cells.Item(8, 2).Value = " September 22 2020"
$ws.Cells.Item(8, 3).Value = "Royals won by 16 runs"
$ws.Cells.Item(8, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(8, 5).Value = "Rajasthan Royals"
$ws.Cells.Item(8, 6).Value = "Ravindra Jadeja "
$ws.Cells.Item(8, 7).Value = "1"
$ws.Cells.Item(8, 8).Value = "2"
$ws.Cells.Item(8, 9).Value = "0"
$ws.Cells.Item(8, 10).Value = "0"
$ws.Cells.Item(8, 11).Value = "50.00"

# Row 9
$ws.Cells.Item(9, 1).Value = " Abu Dhabi"
$ws.Cells.Item(9, 2).Value = " October 19 2020"
$ws.Cells.Item(9, 3).Value = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Cells.Item(9, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(9, 5).Value = "Rajasthan Royals"
$ws.Cells.Item(9, 6).Value = "Ravindra Jadeja "
$ws.Cells.Item(9, 7).Value = "35"
$ws.Cells.Item(9, 8).Value = "30"
$ws.Cells.Item(9, 9).Value = "4"
$ws.Cells.Item(9, 10).Value = "0"
$ws.Cells.Item(9, 11).Value = "116.66"

# Row 10
$ws.Cells.Item(10, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(10, 2).Value = " October 29 2020"
$ws.Cells.Item(10, 3).Value = "Super Kings won by 6 wickets"
$ws.Cells.Item(10, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(10, 5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(10, 6).Value = "Ravindra Jadeja "
$ws.Cells.Item(10, 7).Value = "31"
$ws.Cells.Item(10, 8).Value = "11"
$ws.Cells.Item(10, 9).Value = "2"
$ws.Cells.Item(10, 10).Value = "3"
$ws.Cells.Item(10, 11).Value = "281.81"

# Row 11
$ws.Cells.Item(11, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(11, 2).Value = " October 02 2020"
$ws.Cells.Item(11, 3).Value = "Sunrisers won by 7 runs"
$ws.Cells.Item(11, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(11, 5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(11, 6).Value = "Ravindra Jadeja "
$ws.Cells.Item(11, 7).Value = "50"
$ws.Cells.Item(11, 8).Value = "35"
$ws.Cells.Item(11, 9).Value = "5"
$ws.Cells.Item(11, 10).Value = "2"
$ws.Cells.Item(11, 11).Value = "142.85"

# Row 12
$ws.Cells.Item(12, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(12, 2).Value = " October 10 2020"
$ws.Cells.Item(12, 3).Value = "RCB won by 37 runs"
$ws.Cells.Item(12, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(12, 5).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(12, 6).Value = "Ravindra Jadeja "
$ws.Cells.Item(12, 7).Value = "7"
$ws.Cells.Item(12, 8).Value = "6"
$ws.Cells.Item(12, 9).Value = "1"
$ws.Cells.Item(12, 10).Value = "0"
$ws.Cells.Item(12, 11).Value = "116.66"

Write-Output "done"
